# "update sec debt algorithm"
# The "holding" sheet gains three new columns (LiabilityStartValue,
# InterestRate, DatedDate) inserted right after the existing PosCost_vec
# column (old column H), pushing the former H:M block (Underlying*, Note)
# out to K:P.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("holding")

# Insert three blank columns at H:J — this shifts the old H:M columns
# (and all their data/formatting) to K:P automatically.
$ws.Range("H1:J1").EntireColumn.Insert()

# New header row labels for the inserted columns.
$ws.Range("H1").Value = "LiabilityStartValue"
$ws.Range("I1").Value = "InterestRate"
$ws.Range("J1").Value = "DatedDate"

# Match the column widths used for the new columns.
$ws.Columns.Item(8).ColumnWidth = 16.5
$ws.Columns.Item(9).ColumnWidth = 11.5
$ws.Columns.Item(10).ColumnWidth = 11.5

# Restore the view: scrolled a couple columns in, with H8 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("H8").Select()
